$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target columns keep text storage (avoid numeric auto-conversion)
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "321.47"
$ws.Range("E2").Value = "-3.23%"
$ws.Range("G2").Value = "6"

$ws.Range("D3").Value = "42.95"
$ws.Range("E3").Value = "-6.40%"
$ws.Range("G3").Value = "6"

$ws.Range("D4").Value = "5.184"
$ws.Range("E4").Value = "-8.09%"
$ws.Range("G4").Value = "6"

$ws.Range("D5").Value = "0.08162"
$ws.Range("E5").Value = "-2.53%"
$ws.Range("G5").Value = "6"

$ws.Range("D6").Value = "4.335"
$ws.Range("E6").Value = "-3.32%"
$ws.Range("G6").Value = "6"

$ws.Range("D7").Value = "1.805"
$ws.Range("E7").Value = "-11.57%"
$ws.Range("G7").Value = "6"

$ws.Range("D8").Value = "0.9508"
$ws.Range("E8").Value = "-4.34%"
$ws.Range("G8").Value = "6"

$ws.Range("D9").Value = "0.1114"
$ws.Range("E9").Value = "-3.59%"
$ws.Range("G9").Value = "6"

$ws.Range("D10").Value = "0.1865"
$ws.Range("E10").Value = "-3.66%"
$ws.Range("G10").Value = "6"

$ws.Range("D11").Value = "0.09353"
$ws.Range("E11").Value = "-5.97%"
$ws.Range("G11").Value = "6"

$ws.Range("D12").Value = "0.04635"
$ws.Range("E12").Value = "-1.02%"
$ws.Range("G12").Value = "6"

$ws.Range("D13").Value = "7.416"
$ws.Range("E13").Value = "-28.80%"
$ws.Range("G13").Value = "6"

$ws.Range("D14").Value = "0.1060"
$ws.Range("E14").Value = "0.03%"
$ws.Range("G14").Value = "6"

$ws.Range("D15").Value = "0.001289"
$ws.Range("E15").Value = "1.14%"
$ws.Range("G15").Value = "6"

$ws.Range("D16").Value = "0.005645"
$ws.Range("E16").Value = "-7.21%"
$ws.Range("G16").Value = "6"

$ws.Range("D17").Value = "3.361"
$ws.Range("E17").Value = "-0.46%"
$ws.Range("G17").Value = "6"

$ws.Range("E18").Value = "-2.63%"
$ws.Range("G18").Value = "6"

$ws.Range("D19").Value = "0.3365"
$ws.Range("G19").Value = "6"

$ws.Range("D20").Value = "0.1390"
$ws.Range("E20").Value = "-0.97%"
$ws.Range("G20").Value = "6"

$ws.Range("D21").Value = "0.2628"
$ws.Range("E21").Value = "-0.97%"
$ws.Range("G21").Value = "6"

$ws.Range("D22").Value = "0.04169"
$ws.Range("E22").Value = "-1.20%"
$ws.Range("G22").Value = "6"

$ws.Range("D23").Value = "0.001251"
$ws.Range("E23").Value = "-4.74%"
$ws.Range("G23").Value = "6"

$ws.Range("D24").Value = "0.004297"
$ws.Range("E24").Value = "-7.95%"
$ws.Range("G24").Value = "6"

$ws.Range("E25").Value = "-13.63%"
$ws.Range("G25").Value = "6"

$ws.Range("D26").Value = "0.0002984"
$ws.Range("E26").Value = "-20.49%"
$ws.Range("G26").Value = "6"

$ws.Range("G27").Value = "6"

$ws.Range("G28").Value = "6"

$ws.Range("G29").Value = "6"

$ws.Range("G30").Value = "6"

$ws.Range("G31").Value = "6"

$ws.Range("G32").Value = "6"

$ws.Range("G33").Value = "6"

$ws.Range("G34").Value = "6"

$ws.Range("G35").Value = "6"

$ws.Range("G36").Value = "6"

$ws.Range("G37").Value = "6"

$ws.Range("D38").Value = "0.02674"
$ws.Range("E38").Value = "-4.15%"
$ws.Range("G38").Value = "6"

$ws.Range("D39").Value = "0.05521"
$ws.Range("E39").Value = "-4.12%"
$ws.Range("G39").Value = "6"

$ws.Range("D40").Value = "0.007978"
$ws.Range("E40").Value = "2.25%"
$ws.Range("G40").Value = "6"

$ws.Range("D41").Value = "0.1393"
$ws.Range("E41").Value = "-3.13%"
$ws.Range("G41").Value = "6"

$ws.Range("D42").Value = "0.006558"
$ws.Range("G42").Value = "6"

$ws.Range("D43").Value = "0.002122"
$ws.Range("E43").Value = "5.12%"
$ws.Range("G43").Value = "6"

$ws.Range("D44").Value = "0.008439"
$ws.Range("E44").Value = "-6.75%"
$ws.Range("G44").Value = "6"

$ws.Range("D45").Value = "0.3199"
$ws.Range("E45").Value = "-6.13%"
$ws.Range("G45").Value = "6"

$ws.Range("D46").Value = "0.00006999"
$ws.Range("E46").Value = "-4.54%"
$ws.Range("G46").Value = "6"

$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "-0.18%"
$ws.Range("G47").Value = "6"

$ws.Range("D48").Value = "0.003494"
$ws.Range("E48").Value = "-0.28%"
$ws.Range("G48").Value = "6"

$ws.Range("D49").Value = "0.003536"
$ws.Range("E49").Value = "0.73%"
$ws.Range("G49").Value = "6"

$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "-0.18%"
$ws.Range("G50").Value = "6"

$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "-0.18%"
$ws.Range("G51").Value = "6"
